# Auto-generated edit script: updates cryptos list (D=Price, E=Volume(1h))
# matching commit "Updated cryptos list on Sun Feb 18 05:00:01 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a literal text value to a cell without Excel re-interpreting
# numeric-looking strings (e.g. "352.78") as numbers, and without leaving any
# residual cell-level style (quote-prefix) behind once the text is stored.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "51.569.18"
Set-TextValue $ws.Range("E2") "  -0.87%  "
Set-TextValue $ws.Range("D3") "2.796.63"
Set-TextValue $ws.Range("E3") "  +0.11%  "
Set-TextValue $ws.Range("E4") "  -0.07%  "
Set-TextValue $ws.Range("D5") "352.78"
Set-TextValue $ws.Range("D6") "108.94"
Set-TextValue $ws.Range("E6") "  -0.53%  "
Set-TextValue $ws.Range("D7") "0.553"
Set-TextValue $ws.Range("E7") "  -1.34%  "
Set-TextValue $ws.Range("E8") "  +0.01%  "
Set-TextValue $ws.Range("D9") "0.627"
Set-TextValue $ws.Range("E9") "  +5.57%  "
Set-TextValue $ws.Range("D10") "39.73"
Set-TextValue $ws.Range("E10") "  -1.27%  "
Set-TextValue $ws.Range("E11") "  +0.97%  "
Set-TextValue $ws.Range("D12") "0.0836"
Set-TextValue $ws.Range("E12") "  -1.78%  "
Set-TextValue $ws.Range("D13") "19.98"
Set-TextValue $ws.Range("E13") "  +2.55%  "
Set-TextValue $ws.Range("E14") "  +2.51%  "
Set-TextValue $ws.Range("D15") "3.230.49"
Set-TextValue $ws.Range("E15") "  -0.12%  "
Set-TextValue $ws.Range("D16") "2.788.90"
Set-TextValue $ws.Range("E16") "  +0.00%  "
Set-TextValue $ws.Range("D17") "0.934"
Set-TextValue $ws.Range("E17") "  -1.09%  "
Set-TextValue $ws.Range("D18") "51.584.34"
Set-TextValue $ws.Range("E18") "  -0.73%  "
Set-TextValue $ws.Range("D19") "7.71"
Set-TextValue $ws.Range("E19") "  +2.89%  "
Set-TextValue $ws.Range("E20") "  +1.17%  "
Set-TextValue $ws.Range("D21") "13.33"
Set-TextValue $ws.Range("E21") "  +1.41%  "
Set-TextValue $ws.Range("D22") "0.0₃0970"
Set-TextValue $ws.Range("E22") "  -0.86%  "
Set-TextValue $ws.Range("D23") "70.38"
Set-TextValue $ws.Range("E23") "  +0.19%  "
Set-TextValue $ws.Range("D24") "266.78"
Set-TextValue $ws.Range("E24") "  -1.31%  "
Set-TextValue $ws.Range("D25") "2.76"
Set-TextValue $ws.Range("E25") "  -0.27%  "
Set-TextValue $ws.Range("E26") "  -0.04%  "
Set-TextValue $ws.Range("D27") "25.93"
Set-TextValue $ws.Range("E27") "  -2.33%  "
Set-TextValue $ws.Range("E28") "  +1.74%  "
Set-TextValue $ws.Range("D29") "10.30"
Set-TextValue $ws.Range("E29") "  -0.14%  "
Set-TextValue $ws.Range("D30") "37.12"
Set-TextValue $ws.Range("E30") "  +7.83%  "
Set-TextValue $ws.Range("E31") "  -2.49%  "
Set-TextValue $ws.Range("D32") "6.24"
Set-TextValue $ws.Range("E32") "  +8.54%  "
Set-TextValue $ws.Range("D33") "52.23"
Set-TextValue $ws.Range("E33") "  +0.14%  "
Set-TextValue $ws.Range("D34") "5.66"
Set-TextValue $ws.Range("E34") "  +8.68%  "
Set-TextValue $ws.Range("E35") "  -6.33%  "
Set-TextValue $ws.Range("D36") "0.0851"
Set-TextValue $ws.Range("E36") "  +0.42%  "
Set-TextValue $ws.Range("E37") "  -0.06%  "
Set-TextValue $ws.Range("D38") "18.58"
Set-TextValue $ws.Range("E38") "  -2.49%  "
Set-TextValue $ws.Range("D39") "3.13"
Set-TextValue $ws.Range("E39") "  -2.67%  "
Set-TextValue $ws.Range("E40") "  -1.18%  "
Set-TextValue $ws.Range("E41") "  -0.40%  "
Set-TextValue $ws.Range("D42") "2.49"
Set-TextValue $ws.Range("E42") "  -6.22%  "
Set-TextValue $ws.Range("D43") "119.94"
Set-TextValue $ws.Range("E43") "  +0.50%  "
Set-TextValue $ws.Range("D44") "22.01"
Set-TextValue $ws.Range("E44") "  +0.23%  "
Set-TextValue $ws.Range("E45") "  -2.70%  "
Set-TextValue $ws.Range("D46") "2.130.79"
Set-TextValue $ws.Range("E46") "  +1.94%  "
Set-TextValue $ws.Range("D47") "3.38"
Set-TextValue $ws.Range("E47") "  +3.61%  "
Set-TextValue $ws.Range("E48") "  +6.12%  "
Set-TextValue $ws.Range("D49") "0.224"
Set-TextValue $ws.Range("E49") "  +17.59%  "
Set-TextValue $ws.Range("D50") "0.913"
Set-TextValue $ws.Range("E50") "  -4.51%  "
Set-TextValue $ws.Range("D51") "1.35"
Set-TextValue $ws.Range("E51") "  +9.54%  "
